$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update fecha_final value from 5 to 11, keeping it text (quote-prefixed)
# so the cell keeps its original style/type instead of becoming numeric.
$ws.Range("D2").Value = "'11"

# Move selection to D11 to match the post-edit cursor position
$ws.Range("D11").Select()
